$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder
$ws.Range("H6").Value = 362.1
$ws.Range("I6").Value = 202.33333
$ws.Range("K6").Value = 606.99999
$ws.Range("M6").Value = -494.99999

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2201.303
$ws.Range("J132").Value = 2475.3333
$ws.Range("L132").Value = 7425.999899999999
$ws.Range("N132").Value = -12485.9999

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2568.6
$ws.Range("I137").Value = 1922.7693
$ws.Range("J137").Value = 2950.2273
$ws.Range("K137").Value = 5768.3079
$ws.Range("L137").Value = 8850.6819
$ws.Range("M137").Value = -3218.3079
$ws.Range("N137").Value = -13950.6819

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 4695.7812
$ws.Range("I141").Value = 1621.4517
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 4864.355100000001
$ws.Range("L141").Value = 300000
$ws.Range("M141").Value = 315.6448999999993
$ws.Range("N141").Value = -310360

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 10476.882
$ws.Range("I32").Value = 10474.162
$ws.Range("J32").Value = 10500
$ws.Range("K32").Value = 10474.162
$ws.Range("L32").Value = 10500
$ws.Range("M32").Value = -10187.162
$ws.Range("N32").Value = -11074

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2297.3333
$ws.Range("I61").Value = 2160.16
$ws.Range("K61").Value = 2160.16
$ws.Range("M61").Value = -1948.16

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2297.3333
$ws.Range("I136").Value = 2160.16
$ws.Range("K136").Value = 6480.48
$ws.Range("M136").Value = -3930.48

$ws = $wb.Worksheets.Item("BSM")
# Row 16: Port of Call: Ul'dah
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = $null

# Row 80: Unbreaker
$ws.Range("H80").Value = 1482040.9
$ws.Range("J80").Value = 187.6
$ws.Range("L80").Value = 187.6
$ws.Range("N80").Value = -2183.6

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 1482040.9
$ws.Range("J83").Value = 187.6
$ws.Range("L83").Value = 938
$ws.Range("N83").Value = -10922

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3316.1667
$ws.Range("I105").Value = 3179.6
$ws.Range("K105").Value = 3179.6
$ws.Range("M105").Value = -1432.6

# Row 110: Selective Logging
$ws.Range("H110").Value = 29200
$ws.Range("J110").Value = 29200
$ws.Range("L110").Value = 29200
$ws.Range("N110").Value = -37380

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2584.4412
$ws.Range("I134").Value = 2226.6155
$ws.Range("J134").Value = 3747.375
$ws.Range("K134").Value = 6679.8465
$ws.Range("L134").Value = 11242.125
$ws.Range("M134").Value = -4144.8465
$ws.Range("N134").Value = -16312.125

$ws = $wb.Worksheets.Item("CRP")
# Row 17: Say It with Spears
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1853795
$ws.Range("I58").Value = 2059389
$ws.Range("J58").Value = 3450
$ws.Range("K58").Value = 2059389
$ws.Range("L58").Value = 3450
$ws.Range("M58").Value = -2059186
$ws.Range("N58").Value = -3856

# Row 99: O Pine
$ws.Range("H99").Value = 1610.2222
$ws.Range("I99").Value = 1684.05
$ws.Range("J99").Value = 1399.2858
$ws.Range("K99").Value = 1684.05
$ws.Range("L99").Value = 1399.2858
$ws.Range("M99").Value = -186.05
$ws.Range("N99").Value = -4395.2858

# Row 126: A Better Conductor
$ws.Range("H126").Value = 1610.2222
$ws.Range("I126").Value = 1684.05
$ws.Range("J126").Value = 1399.2858
$ws.Range("K126").Value = 5052.15
$ws.Range("L126").Value = 4197.857400000001
$ws.Range("M126").Value = -2582.15
$ws.Range("N126").Value = -9137.857400000001

# Row 136: Turali Quality
$ws.Range("H136").Value = 1853795
$ws.Range("I136").Value = 2059389
$ws.Range("J136").Value = 3450
$ws.Range("K136").Value = 6178167
$ws.Range("L136").Value = 10350
$ws.Range("M136").Value = -6175617
$ws.Range("N136").Value = -15450

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 3290084.5
$ws.Range("I2").Value = 1275.125
$ws.Range("J2").Value = 5681946
$ws.Range("K2").Value = 7650.75
$ws.Range("L2").Value = 34091676
$ws.Range("M2").Value = -7537.75
$ws.Range("N2").Value = -34091902

# Row 4: In Hot Water
$ws.Range("H4").Value = 334000
$ws.Range("I4").Value = 1000000
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 3000000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -2999888
$ws.Range("N4").Value = -3224

# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 108.35294
$ws.Range("I38").Value = 125.454544
$ws.Range("J38").Value = 77
$ws.Range("K38").Value = 376.363632
$ws.Range("L38").Value = 231
$ws.Range("M38").Value = -29.363632
$ws.Range("N38").Value = -925

# Row 56: Culture Club
$ws.Range("H56").Value = 11069.565
$ws.Range("I56").Value = 11069.565
$ws.Range("K56").Value = 11069.565
$ws.Range("M56").Value = -10539.565

# Row 107: Slippery Service
$ws.Range("H107").Value = 1252.983
$ws.Range("I107").Value = 1236.6471
$ws.Range("J107").Value = 1275.2
$ws.Range("K107").Value = 3709.9413
$ws.Range("L107").Value = 3825.6
$ws.Range("M107").Value = -1789.9413
$ws.Range("N107").Value = -7665.6

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 217986.17
$ws.Range("I113").Value = 250553.08
$ws.Range("J113").Value = 873.5
$ws.Range("K113").Value = 751659.24
$ws.Range("L113").Value = 2620.5
$ws.Range("M113").Value = -749489.24
$ws.Range("N113").Value = -6960.5

$ws = $wb.Worksheets.Item("GSM")
# Row 109: You're My Wonderhall
$ws.Range("H109").Value = 12588
$ws.Range("J109").Value = 12588
$ws.Range("L109").Value = 12588
$ws.Range("N109").Value = -14668

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2294.739
$ws.Range("I122").Value = 2146.3684
$ws.Range("K122").Value = 6439.1052
$ws.Range("M122").Value = -3989.1052

$ws = $wb.Worksheets.Item("LTW")
# Row 13: Throwing Down the Gauntlet
$ws.Range("H13").Value = 29670.666
$ws.Range("I13").Value = 20506
$ws.Range("J13").Value = 48000
$ws.Range("K13").Value = 20506
$ws.Range("L13").Value = 48000
$ws.Range("M13").Value = -20366
$ws.Range("N13").Value = -48280

# Row 14: Quelling Bloody Rumors
$ws.Range("H14").Value = 26500
$ws.Range("I14").Value = 50000
$ws.Range("K14").Value = 50000
$ws.Range("M14").Value = -49828

# Row 74: Overall, We Blend In
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

# Row 77: Eviction Notice (L)
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 3465.5
$ws.Range("I82").Value = 900
$ws.Range("K82").Value = 900
$ws.Range("M82").Value = -539

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 3465.5
$ws.Range("I85").Value = 900
$ws.Range("K85").Value = 900
$ws.Range("M85").Value = 348

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1508
$ws.Range("I93").Value = 1190
$ws.Range("J93").Value = 1826
$ws.Range("K93").Value = 1190
$ws.Range("L93").Value = 1826
$ws.Range("M93").Value = 58
$ws.Range("N93").Value = -4322

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1998.7778
$ws.Range("I132").Value = 1993.6923
$ws.Range("K132").Value = 5981.0769
$ws.Range("M132").Value = -3451.0769
